# Auto-generated script to apply "Updated Vlink test cases" edits
$wb = $excel.ActiveWorkbook

# Update Date (column B) values with new test-run timestamps
$ws = $wb.Worksheets.Item("PayNowCC")
$ws.Range("B2").Value = "Thu Aug 28 22:14:10 EDT 2025"
$ws.Range("B3").Value = "Thu Aug 28 22:15:08 EDT 2025"
$ws.Range("B4").Value = "Thu Aug 28 22:15:34 EDT 2025"
$ws.Range("B5").Value = "Thu Aug 28 22:16:00 EDT 2025"
$ws.Range("B6").Value = "Thu Aug 28 22:16:28 EDT 2025"
$ws.Range("B7").Value = "Thu Aug 28 22:16:53 EDT 2025"
$ws.Range("B8").Value = "Thu Aug 28 22:17:20 EDT 2025"
$ws.Range("B9").Value = "Thu Aug 28 22:17:46 EDT 2025"

$ws = $wb.Worksheets.Item("PayNowCorpSCF")
$ws.Range("B2").Value = "Thu Aug 28 21:51:29 EDT 2025"
$ws.Range("B3").Value = "Thu Aug 28 21:52:11 EDT 2025"
$ws.Range("B4").Value = "Thu Aug 28 21:52:55 EDT 2025"
$ws.Range("B5").Value = "Thu Aug 28 21:53:38 EDT 2025"

$ws = $wb.Worksheets.Item("PayNowCreditSCF")
$ws.Range("B2").Value = "Thu Aug 28 21:57:07 EDT 2025"
$ws.Range("B3").Value = "Thu Aug 28 21:57:48 EDT 2025"
$ws.Range("B4").Value = "Thu Aug 28 21:58:31 EDT 2025"
$ws.Range("B5").Value = "Thu Aug 28 21:59:12 EDT 2025"

$ws = $wb.Worksheets.Item("PayNowPersonalCheckSCF")
$ws.Range("B2").Value = "Thu Aug 28 22:03:09 EDT 2025"
$ws.Range("B3").Value = "Thu Aug 28 22:03:51 EDT 2025"
$ws.Range("B4").Value = "Thu Aug 28 22:04:33 EDT 2025"
$ws.Range("B5").Value = "Thu Aug 28 22:05:15 EDT 2025"

$ws = $wb.Worksheets.Item("PayNowPersonalSavingsSCF")
$ws.Range("B2").Value = "Thu Aug 28 22:08:42 EDT 2025"
$ws.Range("B3").Value = "Thu Aug 28 22:09:24 EDT 2025"
$ws.Range("B4").Value = "Thu Aug 28 22:10:07 EDT 2025"
$ws.Range("B5").Value = "Thu Aug 28 22:10:49 EDT 2025"

$ws = $wb.Worksheets.Item("PayNowCreditDCF")
$ws.Range("B2").Value = "Thu Aug 28 21:54:19 EDT 2025"
$ws.Range("B3").Value = "Thu Aug 28 21:55:01 EDT 2025"
$ws.Range("B4").Value = "Thu Aug 28 21:55:43 EDT 2025"
$ws.Range("B5").Value = "Thu Aug 28 21:56:25 EDT 2025"

$ws = $wb.Worksheets.Item("PayNowCorpDCF")
$ws.Range("B2").Value = "Thu Aug 28 21:48:16 EDT 2025"
$ws.Range("B3").Value = "Thu Aug 28 21:48:58 EDT 2025"
$ws.Range("B4").Value = "Thu Aug 28 21:49:38 EDT 2025"
$ws.Range("B5").Value = "Thu Aug 28 21:50:20 EDT 2025"

$ws = $wb.Worksheets.Item("PayNowPC")
$ws.Range("B2").Value = "Thu Aug 28 21:59:54 EDT 2025"

$ws = $wb.Worksheets.Item("PayNowPersonalCheckDCF")
$ws.Range("B2").Value = "Thu Aug 28 22:05:56 EDT 2025"
$ws.Range("B3").Value = "Thu Aug 28 22:06:38 EDT 2025"
$ws.Range("B4").Value = "Thu Aug 28 22:07:19 EDT 2025"
$ws.Range("B5").Value = "Thu Aug 28 22:08:00 EDT 2025"

$ws = $wb.Worksheets.Item("PayNowPS")
$ws.Range("B2").Value = "Thu Aug 28 22:11:30 EDT 2025"
$ws.Range("B3").Value = "Thu Aug 28 22:11:58 EDT 2025"
$ws.Range("B4").Value = "Thu Aug 28 22:12:23 EDT 2025"
$ws.Range("B5").Value = "Thu Aug 28 22:12:49 EDT 2025"
$ws.Range("B6").Value = "Thu Aug 28 22:13:16 EDT 2025"
$ws.Range("B7").Value = "Thu Aug 28 22:13:43 EDT 2025"

$ws = $wb.Worksheets.Item("OverAndUnderPayCredit")
$ws.Range("B2").Value = "Thu Aug 28 21:42:28 EDT 2025"
$ws.Range("B3").Value = "Thu Aug 28 21:43:25 EDT 2025"
$ws.Range("B4").Value = "Thu Aug 28 21:43:52 EDT 2025"
$ws.Range("B5").Value = "Thu Aug 28 21:44:18 EDT 2025"

$ws = $wb.Worksheets.Item("OverAndUnderPayPC")
$ws.Range("B2").Value = "Thu Aug 28 21:44:45 EDT 2025"
$ws.Range("B3").Value = "Thu Aug 28 21:45:12 EDT 2025"
$ws.Range("B4").Value = "Thu Aug 28 21:45:38 EDT 2025"
$ws.Range("B5").Value = "Thu Aug 28 21:46:04 EDT 2025"

$ws = $wb.Worksheets.Item("OverAndUnderPayPS")
$ws.Range("B2").Value = "Thu Aug 28 21:46:31 EDT 2025"
$ws.Range("B3").Value = "Thu Aug 28 21:46:57 EDT 2025"
$ws.Range("B4").Value = "Thu Aug 28 21:47:22 EDT 2025"
$ws.Range("B5").Value = "Thu Aug 28 21:47:49 EDT 2025"

$ws = $wb.Worksheets.Item("OverAndUnderPayCorp")
$ws.Range("B2").Value = "Thu Aug 28 21:40:04 EDT 2025"
$ws.Range("B3").Value = "Thu Aug 28 21:41:04 EDT 2025"
$ws.Range("B4").Value = "Thu Aug 28 21:41:31 EDT 2025"
$ws.Range("B5").Value = "Thu Aug 28 21:42:01 EDT 2025"

$ws = $wb.Worksheets.Item("PayNowCorp")
$ws.Range("B2").Value = "Thu Aug 28 21:51:02 EDT 2025"

$ws = $wb.Worksheets.Item("CardNotAcceptedErrorCC")
$ws.Range("B2").Value = "Thu Aug 28 22:18:12 EDT 2025"
$ws.Range("B3").Value = "Thu Aug 28 22:18:37 EDT 2025"

# Update Result (column A) from Pass to Fail for specific rows
$wb.Worksheets.Item("PayNowCC").Range("A2").Value = "Fail"
$wb.Worksheets.Item("OverAndUnderPayCredit").Range("A2").Value = "Fail"
$wb.Worksheets.Item("OverAndUnderPayCorp").Range("A2").Value = "Fail"
